$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (prices) but must stay stored
# as literal text, matching how this sheet already stores every value (incl.
# numeric-looking prices) as a string. NumberFormat is applied per cell rather
# than via one combined multi-area Range, since only the first area of a
# multi-area Range reliably receives the format.
$textCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "B13", "C13", "D13", "E13", "B14", "C14", "D14", "E14", "B15", "C15", "D15", "E15", "B16", "C16", "D16", "E16", "B17", "C17", "D17", "E17", "B18", "C18", "D18", "E18", "B19", "C19", "D19", "E19", "B20", "C20", "D20", "E20", "B21", "C21", "D21", "E21", "B22", "C22", "D22", "E22", "B23", "C23", "D23", "E23", "B24", "C24", "D24", "E24", "B25", "C25", "D25", "E25", "B26", "C26", "D26", "E26", "D40", "D41", "E41", "D42", "D43", "D44", "E44", "D45", "E47", "D48", "E48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Apply the updated values row by row ---
# Row 2
$ws.Range("D2").Value = "245.85"

# Row 3
$ws.Range("D3").Value = "24.23"

# Row 4
$ws.Range("D4").Value = "5.362"

# Row 5
$ws.Range("D5").Value = "0.05733"

# Row 6
$ws.Range("D6").Value = "6.501"

# Row 7
$ws.Range("D7").Value = "3.139"

# Row 8
$ws.Range("D8").Value = "0.8164"

# Row 9
$ws.Range("D9").Value = "0.8693"

# Row 10
$ws.Range("D10").Value = "0.1372"

# Row 11
$ws.Range("D11").Value = "0.06992"

# Row 12
$ws.Range("D12").Value = "0.03237"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.02883"
$ws.Range("E13").Value = "12BitrueCoinBTR"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09392"
$ws.Range("E14").Value = "13BitMartTokenBMX"

# Row 15
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "3.747"
$ws.Range("E15").Value = "14MCDexMCB"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001530"
$ws.Range("E16").Value = "15BitForexTokenBF"

# Row 17
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04708"
$ws.Range("E17").Value = "16CoinExTokenCET"

# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005989"
$ws.Range("E18").Value = "17OneONE"

# Row 19
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "0.006155"
$ws.Range("E19").Value = "18TigerCashTCH"

# Row 20
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "0.001245"
$ws.Range("E20").Value = "19BitKanKAN"

# Row 21
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "0.004780"
$ws.Range("E21").Value = "20HotbitTokenHTB"

# Row 22
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "0.00006798"
$ws.Range("E22").Value = "21NitroExNTX"

# Row 23
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "3.529"
$ws.Range("E23").Value = "22LEOLEO"

# Row 24
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "2.148"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# Row 25
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").Value = "0.3155"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"

# Row 26
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").Value = "0.1331"
$ws.Range("E26").Value = "25ProBitTokenPROB"

# Row 40
$ws.Range("D40").Value = "0.03705"

# Row 41
$ws.Range("D41").Value = "0.006408"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

# Row 42
$ws.Range("D42").Value = "0.1055"

# Row 43
$ws.Range("D43").Value = "0.002211"

# Row 44
$ws.Range("D44").Value = "0.008637"
$ws.Range("E44").Value = "43LocalTradersLCT"

# Row 45
$ws.Range("D45").Value = "0.00005490"

# Row 47
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

# Row 48
$ws.Range("D48").Value = "0.002564"
$ws.Range("E48").Value = "47BOLOBOLO"

# Restore the default (unstyled) cell style now that the text has been written,
# so we do not leave a stray custom number format behind on these cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
